# Update the "bmwModels" worksheet: the header row is trimmed down to a
# single "Model" column and the first data row now holds the first model
# value ("X1") directly in column A. Column B is cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bmwModels")

# Header row: A1 becomes "Model", B1 is cleared.
$ws.Range("A1").Value = "Model"
$ws.Range("B1").ClearContents()

# First data row: A2 becomes "X1", B2 is cleared.
$ws.Range("A2").Value = "X1"
$ws.Range("B2").ClearContents()

# Update the active selection to B1, matching the saved workbook state.
$ws.Activate()
$ws.Range("B1").Select()
